$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: header "IC_Br" + values for rows 2-25
$ws.Range("E1").Value = "IC_Br"

$icBrValues = @{
    2 = 10
    3 = -19
    4 = 14
    5 = 15
    6 = 12
    7 = 17
    8 = 7
    9 = 12
    10 = -22
    11 = 28
    12 = 28
    13 = -8
    14 = 0
    15 = -9
    16 = -3
    17 = -16
    18 = 11
    19 = 1
    20 = -3
    21 = 0
    22 = 2
    23 = 37
    24 = 6
    25 = -6
}

foreach ($row in $icBrValues.Keys | Sort-Object) {
    $ws.Cells.Item($row, 5).Value = $icBrValues[$row]
}

# View changes: zoom to 125% and move the active selection to C11
$excel.ActiveWindow.Zoom = 125
$ws.Range("C11").Select()
